# Fruta / hortaliza, semanal
# Insert two new weekly price rows (at the top of the data block, row 514)
# for "Femacal de La Calera" / Mango, pushing the existing rows 514-585
# down to 516-587.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 514 (existing data shifts down by 2).
$ws.Range("A514:A515").EntireRow.Insert()

# New row 514
$ws.Range("A514").Value = 3
$ws.Range("B514").Value = "Femacal de La Calera"
$ws.Range("C514").Value = "Coquimbo"
$ws.Range("D514").Value = 44984
$ws.Range("E514").Value = 5
$ws.Range("F514").Value = "Fruta"
$ws.Range("G514").Value = 100108
$ws.Range("H514").Value = "Tropicales y subtropicales"
$ws.Range("I514").Value = 100108002
$ws.Range("J514").Value = "Mango"
$ws.Range("K514").Value = "Sin especificar"
$ws.Range("L514").Value = "Primera"
$ws.Range("M514").Value = 228
$ws.Range("N514").Value = 7000
$ws.Range("O514").Value = 7000
$ws.Range("P514").Value = 7000
$ws.Range("Q514").Value = "$/bandeja 4 kilos"
$ws.Range("R514").Value = "Perú"
$ws.Range("S514").Value = 1750
$ws.Range("T514").Value = 4

# New row 515
$ws.Range("A515").Value = 3
$ws.Range("B515").Value = "Femacal de La Calera"
$ws.Range("C515").Value = "Coquimbo"
$ws.Range("D515").Value = 44984
$ws.Range("E515").Value = 5
$ws.Range("F515").Value = "Fruta"
$ws.Range("G515").Value = 100108
$ws.Range("H515").Value = "Tropicales y subtropicales"
$ws.Range("I515").Value = 100108002
$ws.Range("J515").Value = "Mango"
$ws.Range("K515").Value = "Sin especificar"
$ws.Range("L515").Value = "Segunda"
$ws.Range("M515").Value = 228
$ws.Range("N515").Value = 7000
$ws.Range("O515").Value = 7000
$ws.Range("P515").Value = 7000
$ws.Range("Q515").Value = "$/bandeja 4 kilos"
$ws.Range("R515").Value = "Perú"
$ws.Range("S515").Value = 1750
$ws.Range("T515").Value = 4
